# Progress_Tracking_sheet_9d.xlsx edit
# Updates user-story texts, effort/priority numbers and the "Done?" flags on
# Sheet1, and refreshes the sheet view (scroll/zoom/selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- User story text (column C) updates -----------------------------------
# NOTE: these five assignments introduce brand-new shared strings, so they
# are ordered to match the final shared-string table layout (row 5, 15, 16,
# 21, then 17).
$ws.Range("C5").Value = "As a website visitor, I would like to receive a confirmation that I have booked, to ensure my booking took place."
$ws.Range("C15").Value = "As a tourist, I would like to be able to search for tours that do pickup, so I can plan my day efficiently."
$ws.Range("C16").Value = "As a marketing manager, I would like to be able to make a new trip, so it can be viewed on the site"
$ws.Range("C21").Value = "As a traveler, I would like to search for tours depending on their difficulty, so I can filter out tours I am not fit enough for"
$ws.Range("C17").Value = "As a marketing manager, I want to be able to accept all pending reviews, to save time when I know reviews are not spam"

# --- Actual Team Effort (column F) + Estimated Team Effort (column E) -----
$ws.Range("F4").Value = 3

$ws.Range("F7").Value = 1
$ws.Range("F8").Value = 5
$ws.Range("F9").Value = 2

$ws.Range("F11").Value = 1
$ws.Range("F12").Value = 1
$ws.Range("F13").Value = 1
$ws.Range("F14").Value = 1

$ws.Range("F15").Value = 1

$ws.Range("E16").Value = 5
$ws.Range("F16").Value = 2

$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 1

$ws.Range("E21").Value = 1
$ws.Range("F21").Value = 1

# --- Done? flags (column G) ------------------------------------------------
$ws.Range("G5").Value = "YES"
$ws.Range("G15").Value = "YES"
$ws.Range("G16").Value = "YES"
$ws.Range("G17").Value = "YES"
$ws.Range("G21").Value = "YES"

# --- Sheet view: scroll position, zoom, selection --------------------------
$ws.Select()
$win = $excel.ActiveWindow
$win.ScrollColumn = 3
$win.ScrollRow = 3
$win.Zoom = 90
$ws.Range("C23").Select()
